$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.114.65"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").Value = "2.500.89"
$ws.Range("E3").Value = "  -0.85%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'319.70"
$ws.Range("E5").Value = "  -1.17%  "

$ws.Range("D6").Value = "'106.40"
$ws.Range("E6").Value = "  -2.79%  "

$ws.Range("D7").Value = "'0.523"
$ws.Range("E7").Value = "  -0.53%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -3.87%  "

$ws.Range("D10").Value = "'38.82"
$ws.Range("E10").Value = "  -3.72%  "

$ws.Range("D11").Value = "'20.03"
$ws.Range("E11").Value = "  +1.99%  "

$ws.Range("E12").Value = "  -2.13%  "

$ws.Range("E13").Value = "  -0.66%  "

$ws.Range("E14").Value = "  -2.15%  "

$ws.Range("D15").Value = "2.893.47"
$ws.Range("E15").Value = "  -0.49%  "

$ws.Range("D16").Value = "2.501.21"
$ws.Range("E16").Value = "  -0.61%  "

$ws.Range("D17").Value = "'0.833"
$ws.Range("E17").Value = "  -2.59%  "

$ws.Range("D18").Value = "48.009.99"
$ws.Range("E18").Value = "  -0.41%  "

$ws.Range("D19").Value = "'12.97"
$ws.Range("E19").Value = "  -3.50%  "

$ws.Range("E20").Value = "  +8.15%  "

$ws.Range("D21").Value = "'6.64"
$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("E22").Value = "  -1.16%  "

$ws.Range("D23").Value = "'71.21"
$ws.Range("E23").Value = "  -1.75%  "

$ws.Range("D24").Value = "'272.51"
$ws.Range("E24").Value = "  +1.85%  "

$ws.Range("D25").Value = "'2.52"
$ws.Range("E25").Value = "  -2.04%  "

$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("D27").Value = "'25.80"
$ws.Range("E27").Value = "  -1.62%  "

$ws.Range("E28").Value = "  +10.16%  "

$ws.Range("E29").Value = "  -1.30%  "

$ws.Range("D30").Value = "'9.73"
$ws.Range("E30").Value = "  -4.64%  "

$ws.Range("D31").Value = "'35.05"
$ws.Range("E31").Value = "  -0.90%  "

$ws.Range("D32").Value = "'49.33"
$ws.Range("E32").Value = "  -1.24%  "

$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("D34").Value = "'19.12"
$ws.Range("E34").Value = "  -4.65%  "

$ws.Range("D35").Value = "'5.29"
$ws.Range("E35").Value = "  -2.14%  "

$ws.Range("D36").Value = "'0.0776"
$ws.Range("E36").Value = "  -1.52%  "

$ws.Range("E37").Value = "  -2.68%  "

$ws.Range("D38").Value = "'4.58"
$ws.Range("E38").Value = "  -2.82%  "

$ws.Range("D39").Value = "'2.87"
$ws.Range("E39").Value = "  -3.53%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "'121.70"
$ws.Range("E40").Value = "  +2.40%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "'0.111"
$ws.Range("E41").Value = "  -1.53%  "

$ws.Range("D42").Value = "'22.04"
$ws.Range("E42").Value = "  -1.07%  "

$ws.Range("D43").Value = "'2.20"
$ws.Range("E43").Value = "  +1.18%  "

$ws.Range("E44").Value = "  +2.02%  "

$ws.Range("D45").Value = "2.007.29"
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("D46").Value = "'3.17"
$ws.Range("E46").Value = "  +1.85%  "

$ws.Range("E47").Value = "  +2.16%  "

$ws.Range("D49").Value = "'8.94"
$ws.Range("E49").Value = "  -1.73%  "

$ws.Range("E50").Value = "  -1.66%  "

$ws.Range("D51").Value = "'78.92"
$ws.Range("E51").Value = "  -1.78%  "
